# Generate Report for Archive
#
# Two files (749b72bf-41b4-4a05-bcd0-4eeabaf9f8dd and
# a5460bc7-220b-4afd-a261-56e41741d2d2) have moved out of "Ready for
# handoff" and are now "In Translation" for both locales. Additionally,
# d4dc8520-4bb8-4b24-87fb-2b204206540d has moved to "In Translation" for
# the de-de locale only. Update the per-locale Status columns as well as
# the locale-status columns on the Overview summary sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "In Translation"
$overview.Range("C3").Value = "In Translation"
$overview.Range("B4").Value = "In Translation"
$overview.Range("C4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"
$dede.Range("C5").Value = "In Translation"
